$d = $word.ActiveDocument

$pairs = @(
    @("2024-04-06 Saturday", "2024-04-07 Sunday"),
    @("63×42=2646", "25×46=1150"),
    @("38×47=1786", "69×48=3312"),
    @("67×57=3819", "72×23=1656"),
    @("98×27=2646", "29×20=580"),
    @("97×98=9506", "18×47=846"),
    @("24×54=1296", "19×49=931"),
    @("40×40=1600", "99×26=2574"),
    @("23×53=1219", "42×87=3654"),
    @("87×72=6264", "56×74=4144"),
    @("23×86=1978", "26×95=2470"),
    @("62×34=2108", "64×98=6272"),
    @("51×15=765", "62×32=1984"),
    @("54×48=2592", "37×54=1998"),
    @("29×84=2436", "76×75=5700"),
    @("63×55=3465", "90×27=2430"),
    @("36×89=3204", "37×42=1554"),
    @("76×73=5548", "30×38=1140"),
    @("70×47=3290", "61×79=4819"),
    @("89×32=2848", "14×52=728"),
    @("73×71=5183", "57×15=855"),
    @("47×83=3901", "83×35=2905"),
    @("38×75=2850", "88×72=6336"),
    @("69×14=966", "18×78=1404"),
    @("35×59=2065", "42×18=756"),
    @("76×58=4408", "13×36=468")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
